$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 205
$ws.Range("I2").Value = 205
$ws.Range("K2").Value = 205
$ws.Range("M2").Value = -92
$ws.Range("H9").Value = 148.72728
$ws.Range("I9").Value = 81.30768999999999
$ws.Range("J9").Value = 246.11111
$ws.Range("K9").Value = 81.30768999999999
$ws.Range("L9").Value = 246.11111
$ws.Range("M9").Value = 87.69231000000001
$ws.Range("N9").Value = -584.1111100000001
$ws.Range("H12").Value = 1122.8
$ws.Range("I12").Value = 175.57143
$ws.Range("J12").Value = 3333
$ws.Range("K12").Value = 175.57143
$ws.Range("L12").Value = 3333
$ws.Range("M12").Value = -5.571429999999992
$ws.Range("N12").Value = -3673
$ws.Range("H43").Value = 4131.375
$ws.Range("J43").Value = 4300.1665
$ws.Range("L43").Value = 4300.1665
$ws.Range("N43").Value = -4438.1665
$ws.Range("H52").Value = 4161.8
$ws.Range("I52").Value = 4352.25
$ws.Range("K52").Value = 13056.75
$ws.Range("M52").Value = -12896.75
$ws.Range("H69").Value = 63753
$ws.Range("J69").Value = 62500
$ws.Range("L69").Value = 187500
$ws.Range("N69").Value = -189248
$ws.Range("H70").Value = 6035.524
$ws.Range("I70").Value = 3081.1428
$ws.Range("J70").Value = 7512.7144
$ws.Range("K70").Value = 9243.428400000001
$ws.Range("L70").Value = 22538.1432
$ws.Range("M70").Value = -8973.428400000001
$ws.Range("N70").Value = -23078.1432
$ws.Range("H72").Value = 63753
$ws.Range("J72").Value = 62500
$ws.Range("L72").Value = 562500
$ws.Range("N72").Value = -571236
$ws.Range("H73").Value = 6035.524
$ws.Range("I73").Value = 3081.1428
$ws.Range("J73").Value = 7512.7144
$ws.Range("K73").Value = 9243.428400000001
$ws.Range("L73").Value = 22538.1432
$ws.Range("M73").Value = -8307.428400000001
$ws.Range("N73").Value = -24410.1432
$ws.Range("H76").Value = 9073.666999999999
$ws.Range("J76").Value = 11110.5
$ws.Range("L76").Value = 11110.5
$ws.Range("N76").Value = -11740.5
$ws.Range("H79").Value = 9073.666999999999
$ws.Range("J79").Value = 11110.5
$ws.Range("L79").Value = 11110.5
$ws.Range("N79").Value = -13294.5
$ws.Range("H82").Value = 1212.3889
$ws.Range("I82").Value = 1048.4117
$ws.Range("K82").Value = 3145.2351
$ws.Range("M82").Value = -2739.2351
$ws.Range("H85").Value = 1212.3889
$ws.Range("I85").Value = 1048.4117
$ws.Range("K85").Value = 3145.2351
$ws.Range("M85").Value = -1741.2351
$ws.Range("H88").Value = 4905.0713
$ws.Range("I88").Value = 7398.4
$ws.Range("J88").Value = 3519.889
$ws.Range("K88").Value = 7398.4
$ws.Range("L88").Value = 3519.889
$ws.Range("M88").Value = -6992.4
$ws.Range("N88").Value = -4331.889
$ws.Range("H91").Value = 4905.0713
$ws.Range("I91").Value = 7398.4
$ws.Range("J91").Value = 3519.889
$ws.Range("K91").Value = 7398.4
$ws.Range("L91").Value = 3519.889
$ws.Range("M91").Value = -5994.4
$ws.Range("N91").Value = -6327.889
$ws.Range("H138").Value = 4833.1724
$ws.Range("I138").Value = 4284.25
$ws.Range("J138").Value = 5042.2856
$ws.Range("K138").Value = 12852.75
$ws.Range("L138").Value = 15126.8568
$ws.Range("M138").Value = -7712.75
$ws.Range("N138").Value = -25406.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 12857.143
$ws.Range("J23").Value = 12857.143
$ws.Range("L23").Value = 12857.143
$ws.Range("N23").Value = -13375.143
$ws.Range("H32").Value = 3389.3289
$ws.Range("I32").Value = 3431.5571
$ws.Range("J32").Value = 2404
$ws.Range("K32").Value = 3431.5571
$ws.Range("L32").Value = 2404
$ws.Range("M32").Value = -3144.5571
$ws.Range("N32").Value = -2978
$ws.Range("H74").Value = 4667.7
$ws.Range("I74").Value = 4667.7
$ws.Range("K74").Value = 4667.7
$ws.Range("M74").Value = -3793.7
$ws.Range("H77").Value = 4667.7
$ws.Range("I77").Value = 4667.7
$ws.Range("K77").Value = 23338.5
$ws.Range("M77").Value = -18970.5
$ws.Range("H110").Value = 3713.35
$ws.Range("I110").Value = 3018.625
$ws.Range("K110").Value = 3018.625
$ws.Range("M110").Value = -973.625
$ws.Range("H132").Value = 4109.641
$ws.Range("I132").Value = 2294.9614
$ws.Range("K132").Value = 6884.8842
$ws.Range("M132").Value = -4354.8842

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 49582.57
$ws.Range("J86").Value = 2426.5715
$ws.Range("L86").Value = 2426.5715
$ws.Range("N86").Value = -4672.5715
$ws.Range("H89").Value = 49582.57
$ws.Range("J89").Value = 2426.5715
$ws.Range("L89").Value = 12132.8575
$ws.Range("N89").Value = -23364.8575
$ws.Range("H94").Value = 2755.258
$ws.Range("I94").Value = 2380.2964
$ws.Range("K94").Value = 2380.2964
$ws.Range("M94").Value = -1929.2964
$ws.Range("H99").Value = 34831
$ws.Range("I99").Value = 52404.25
$ws.Range("J99").Value = 11400
$ws.Range("K99").Value = 52404.25
$ws.Range("L99").Value = 11400
$ws.Range("M99").Value = -50906.25
$ws.Range("N99").Value = -14396

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5107.4053
$ws.Range("I31").Value = 5929.25
$ws.Range("J31").Value = 3590.1538
$ws.Range("K31").Value = 5929.25
$ws.Range("L31").Value = 3590.1538
$ws.Range("M31").Value = -5634.25
$ws.Range("N31").Value = -4180.1538
$ws.Range("H34").Value = 5107.4053
$ws.Range("I34").Value = 5929.25
$ws.Range("J34").Value = 3590.1538
$ws.Range("K34").Value = 5929.25
$ws.Range("L34").Value = 3590.1538
$ws.Range("M34").Value = -5727.25
$ws.Range("N34").Value = -3994.1538
$ws.Range("H132").Value = 3064.2856
$ws.Range("I132").Value = 2761.6155
$ws.Range("K132").Value = 8284.8465
$ws.Range("M132").Value = -5754.8465

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 820.1
$ws.Range("J86").Value = 831.6667
$ws.Range("L86").Value = 2495.0001
$ws.Range("N86").Value = -4867.0001
$ws.Range("H89").Value = 820.1
$ws.Range("J89").Value = 831.6667
$ws.Range("L89").Value = 7485.0003
$ws.Range("N89").Value = -19341.0003
$ws.Range("H113").Value = 864.3
$ws.Range("I113").Value = 1543.3334
$ws.Range("J113").Value = 573.2857
$ws.Range("K113").Value = 4630.0002
$ws.Range("L113").Value = 1719.8571
$ws.Range("M113").Value = -2460.0002
$ws.Range("N113").Value = -6059.8571
$ws.Range("H116").Value = 4400
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").Value = ""
$ws.Range("H133").Value = 3500
$ws.Range("I133").Value = 3500
$ws.Range("K133").Value = 10500
$ws.Range("M133").Value = -5440

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 9005499
$ws.Range("I20").Value = 11251000
$ws.Range("J20").Value = 23495
$ws.Range("K20").Value = 11251000
$ws.Range("L20").Value = 23495
$ws.Range("M20").Value = -11250755
$ws.Range("N20").Value = -23985
$ws.Range("H24").Value = 6436713
$ws.Range("J24").Value = 14247.5
$ws.Range("L24").Value = 14247.5
$ws.Range("N24").Value = -14593.5
$ws.Range("H52").Value = 26713.143
$ws.Range("J52").Value = 26713.143
$ws.Range("L52").Value = 26713.143
$ws.Range("N52").Value = -27231.143
$ws.Range("H80").Value = 72296.44
$ws.Range("I80").Value = 224482.2
$ws.Range("J80").Value = 3121.0908
$ws.Range("K80").Value = 224482.2
$ws.Range("L80").Value = 3121.0908
$ws.Range("M80").Value = -223484.2
$ws.Range("N80").Value = -5117.0908
$ws.Range("H83").Value = 72296.44
$ws.Range("I83").Value = 224482.2
$ws.Range("J83").Value = 3121.0908
$ws.Range("K83").Value = 1122411
$ws.Range("L83").Value = 15605.454
$ws.Range("M83").Value = -1117419
$ws.Range("N83").Value = -25589.454
$ws.Range("H95").Value = 26289
$ws.Range("J95").Value = 26289
$ws.Range("L95").Value = 26289
$ws.Range("N95").Value = -31781
$ws.Range("H98").Value = 22421.5
$ws.Range("J98").Value = 22421.5
$ws.Range("L98").Value = 22421.5
$ws.Range("N98").Value = -28411.5
$ws.Range("H132").Value = 5515.684
$ws.Range("I132").Value = 4814.531
$ws.Range("K132").Value = 14443.593
$ws.Range("M132").Value = -11913.593
$ws.Range("H135").Value = 64811.875
$ws.Range("J135").Value = 64811.875
$ws.Range("L135").Value = 64811.875
$ws.Range("N135").Value = -74951.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 978.2
$ws.Range("I22").Value = 1000.5
$ws.Range("J22").Value = 963.3333
$ws.Range("K22").Value = 1000.5
$ws.Range("L22").Value = 963.3333
$ws.Range("M22").Value = -705.5
$ws.Range("N22").Value = -1553.3333
$ws.Range("H27").Value = 978.2
$ws.Range("I27").Value = 1000.5
$ws.Range("J27").Value = 963.3333
$ws.Range("K27").Value = 1000.5
$ws.Range("L27").Value = 963.3333
$ws.Range("M27").Value = -893.5
$ws.Range("N27").Value = -1177.3333
$ws.Range("H55").Value = 1087.2222
$ws.Range("I55").Value = 540.7143
$ws.Range("K55").Value = 540.7143
$ws.Range("M55").Value = -367.7143
$ws.Range("H122").Value = 6752
$ws.Range("I122").Value = 6752
$ws.Range("K122").Value = 20256
$ws.Range("M122").Value = -17806
$ws.Range("H136").Value = 2661.325
$ws.Range("J136").Value = 3785.2
$ws.Range("L136").Value = 11355.6
$ws.Range("N136").Value = -16455.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 4001
$ws.Range("I31").Value = 4001
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 4001
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -3653
$ws.Range("N31").Value = ""
$ws.Range("H96").Value = 61344
$ws.Range("I96").Value = 103819.4
$ws.Range("K96").Value = 103819.4
$ws.Range("M96").Value = -102446.4
